$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 307-308; this pushes the existing rows
# 307..410 down to 309..412 and grows the used range to A1:R412,
# matching the "Fruta / hortaliza, semanal" weekly-roll edit.
$ws.Range("307:308").EntireRow.Insert()

# Row 307: new weekly record (Primera)
$ws.Cells.Item(307, 1).Value = 3
$ws.Cells.Item(307, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(307, 3).Value = "Coquimbo"
$ws.Cells.Item(307, 4).Value = 44524
$ws.Cells.Item(307, 5).Value = 5
$ws.Cells.Item(307, 6).Value = 100112006
$ws.Cells.Item(307, 7).Value = "Repollo"
$ws.Cells.Item(307, 8).Value = "Crespo record"
$ws.Cells.Item(307, 9).Value = "Primera"
$ws.Cells.Item(307, 10).Value = 1600
$ws.Cells.Item(307, 11).Value = 600
$ws.Cells.Item(307, 12).Value = 600
$ws.Cells.Item(307, 13).Value = 600
$ws.Cells.Item(307, 14).Value = "`$/unidad"
$ws.Cells.Item(307, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(307, 16).Value = 600
$ws.Cells.Item(307, 17).Value = 1
$ws.Cells.Item(307, 18).Value = "Hortaliza"

# Row 308: new weekly record (Segunda)
$ws.Cells.Item(308, 1).Value = 3
$ws.Cells.Item(308, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(308, 3).Value = "Coquimbo"
$ws.Cells.Item(308, 4).Value = 44524
$ws.Cells.Item(308, 5).Value = 5
$ws.Cells.Item(308, 6).Value = 100112006
$ws.Cells.Item(308, 7).Value = "Repollo"
$ws.Cells.Item(308, 8).Value = "Crespo record"
$ws.Cells.Item(308, 9).Value = "Segunda"
$ws.Cells.Item(308, 10).Value = 800
$ws.Cells.Item(308, 11).Value = 500
$ws.Cells.Item(308, 12).Value = 500
$ws.Cells.Item(308, 13).Value = 500
$ws.Cells.Item(308, 14).Value = "`$/unidad"
$ws.Cells.Item(308, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(308, 16).Value = 500
$ws.Cells.Item(308, 17).Value = 1
$ws.Cells.Item(308, 18).Value = "Hortaliza"
